# GOM_GuessOMeter — add the August 2025 consumption history sheet.
# 1) Rename the existing (June) history sheet.
# 2) Insert a fresh sheet right after it for August, carrying over the
#    same layout/formulas as June but with the new trip data.
# 3) Recreate the summary block (sum / average / estimated range / GOM
#    comparison) and the "good match" conclusion line for August.

$wb = $excel.ActiveWorkbook

$wsJune = $wb.Worksheets.Item(1)
$wsJune.Name = "historyJune2025"
$wsJune.Range("A1").Select()

$wsAug = $wb.Worksheets.Add($null, $wsJune)
$wsAug.Name = "historyAugust2025"

# -- headers (row 5) -------------------------------------------------
$wsAug.Cells.Item(5,1).Value = "index"
$wsAug.Cells.Item(5,3).Value = "km"
$wsAug.Cells.Item(5,4).Value = "kWh/100km"
$wsAug.Cells.Item(5,5).Value = "kWh"

# -- trip log (rows 6-35): index, km, kWh/100km, computed kWh --------
$wsAug.Cells.Item(6,1).Value = 1
$wsAug.Cells.Item(6,3).Value = 7
$wsAug.Cells.Item(6,4).Value = 11.7
$wsAug.Cells.Item(6,5).Formula = "=D6/100*C6"
$wsAug.Cells.Item(7,1).Value = 2
$wsAug.Cells.Item(7,3).Value = 794
$wsAug.Cells.Item(7,4).Value = 11.7
$wsAug.Cells.Item(7,5).Formula = "=D7/100*C7"
$wsAug.Cells.Item(8,1).Value = 3
$wsAug.Cells.Item(8,3).Value = 65
$wsAug.Cells.Item(8,4).Value = 10.1
$wsAug.Cells.Item(8,5).Formula = "=D8/100*C8"
$wsAug.Cells.Item(9,1).Value = 4
$wsAug.Cells.Item(9,3).Value = 66
$wsAug.Cells.Item(9,4).Value = 9.4
$wsAug.Cells.Item(9,5).Formula = "=D9/100*C9"
$wsAug.Cells.Item(10,1).Value = 5
$wsAug.Cells.Item(10,3).Value = 44
$wsAug.Cells.Item(10,4).Value = 10.8
$wsAug.Cells.Item(10,5).Formula = "=D10/100*C10"
$wsAug.Cells.Item(11,1).Value = 6
$wsAug.Cells.Item(11,3).Value = 38
$wsAug.Cells.Item(11,4).Value = 8.6
$wsAug.Cells.Item(11,5).Formula = "=D11/100*C11"
$wsAug.Cells.Item(12,1).Value = 7
$wsAug.Cells.Item(12,3).Value = 63
$wsAug.Cells.Item(12,4).Value = 9.8
$wsAug.Cells.Item(12,5).Formula = "=D12/100*C12"
$wsAug.Cells.Item(13,1).Value = 8
$wsAug.Cells.Item(13,3).Value = 66
$wsAug.Cells.Item(13,4).Value = 10.4
$wsAug.Cells.Item(13,5).Formula = "=D13/100*C13"
$wsAug.Cells.Item(14,1).Value = 9
$wsAug.Cells.Item(14,3).Value = 266
$wsAug.Cells.Item(14,4).Value = 11.4
$wsAug.Cells.Item(14,5).Formula = "=D14/100*C14"
$wsAug.Cells.Item(15,1).Value = 10
$wsAug.Cells.Item(15,3).Value = 81
$wsAug.Cells.Item(15,4).Value = 10.4
$wsAug.Cells.Item(15,5).Formula = "=D15/100*C15"
$wsAug.Cells.Item(16,1).Value = 11
$wsAug.Cells.Item(16,3).Value = 672
$wsAug.Cells.Item(16,4).Value = 11.8
$wsAug.Cells.Item(16,5).Formula = "=D16/100*C16"
$wsAug.Cells.Item(17,1).Value = 12
$wsAug.Cells.Item(17,3).Value = 973
$wsAug.Cells.Item(17,4).Value = 13.4
$wsAug.Cells.Item(17,5).Formula = "=D17/100*C17"
$wsAug.Cells.Item(18,1).Value = 13
$wsAug.Cells.Item(18,3).Value = 2136
$wsAug.Cells.Item(18,4).Value = 11.5
$wsAug.Cells.Item(18,5).Formula = "=D18/100*C18"
$wsAug.Cells.Item(19,1).Value = 14
$wsAug.Cells.Item(19,3).Value = 84
$wsAug.Cells.Item(19,4).Value = 11.1
$wsAug.Cells.Item(19,5).Formula = "=D19/100*C19"
$wsAug.Cells.Item(20,1).Value = 15
$wsAug.Cells.Item(20,3).Value = 1534
$wsAug.Cells.Item(20,4).Value = 11.9
$wsAug.Cells.Item(20,5).Formula = "=D20/100*C20"
$wsAug.Cells.Item(21,1).Value = 16
$wsAug.Cells.Item(21,3).Value = 1088
$wsAug.Cells.Item(21,4).Value = 13.7
$wsAug.Cells.Item(21,5).Formula = "=D21/100*C21"
$wsAug.Cells.Item(22,1).Value = 17
$wsAug.Cells.Item(22,3).Value = 1066
$wsAug.Cells.Item(22,4).Value = 13.6
$wsAug.Cells.Item(22,5).Formula = "=D22/100*C22"
$wsAug.Cells.Item(23,1).Value = 18
$wsAug.Cells.Item(23,3).Value = 7
$wsAug.Cells.Item(23,4).Value = 13.6
$wsAug.Cells.Item(23,5).Formula = "=D23/100*C23"
$wsAug.Cells.Item(24,1).Value = 19
$wsAug.Cells.Item(24,3).Value = 359
$wsAug.Cells.Item(24,4).Value = 12
$wsAug.Cells.Item(24,5).Formula = "=D24/100*C24"
$wsAug.Cells.Item(25,1).Value = 20
$wsAug.Cells.Item(25,3).Value = 58
$wsAug.Cells.Item(25,4).Value = 9.6
$wsAug.Cells.Item(25,5).Formula = "=D25/100*C25"
$wsAug.Cells.Item(26,1).Value = 21
$wsAug.Cells.Item(26,3).Value = 296
$wsAug.Cells.Item(26,4).Value = 12.8
$wsAug.Cells.Item(26,5).Formula = "=D26/100*C26"
$wsAug.Cells.Item(27,1).Value = 22
$wsAug.Cells.Item(27,3).Value = 79
$wsAug.Cells.Item(27,4).Value = 12
$wsAug.Cells.Item(27,5).Formula = "=D27/100*C27"
$wsAug.Cells.Item(28,1).Value = 23
$wsAug.Cells.Item(28,3).Value = 17
$wsAug.Cells.Item(28,4).Value = 12.7
$wsAug.Cells.Item(28,5).Formula = "=D28/100*C28"
$wsAug.Cells.Item(29,1).Value = 24
$wsAug.Cells.Item(29,3).Value = 16
$wsAug.Cells.Item(29,4).Value = 12.2
$wsAug.Cells.Item(29,5).Formula = "=D29/100*C29"
$wsAug.Cells.Item(30,1).Value = 25
$wsAug.Cells.Item(30,3).Value = 79
$wsAug.Cells.Item(30,4).Value = 11.8
$wsAug.Cells.Item(30,5).Formula = "=D30/100*C30"
$wsAug.Cells.Item(31,1).Value = 26
$wsAug.Cells.Item(31,3).Value = 196
$wsAug.Cells.Item(31,4).Value = 11.8
$wsAug.Cells.Item(31,5).Formula = "=D31/100*C31"
$wsAug.Cells.Item(32,1).Value = 27
$wsAug.Cells.Item(32,3).Value = 8
$wsAug.Cells.Item(32,4).Value = 9.7
$wsAug.Cells.Item(32,5).Formula = "=D32/100*C32"
$wsAug.Cells.Item(33,1).Value = 28
$wsAug.Cells.Item(33,3).Value = 4
$wsAug.Cells.Item(33,4).Value = 15.5
$wsAug.Cells.Item(33,5).Formula = "=D33/100*C33"
$wsAug.Cells.Item(34,1).Value = 29
$wsAug.Cells.Item(34,3).Value = 40
$wsAug.Cells.Item(34,4).Value = 12
$wsAug.Cells.Item(34,5).Formula = "=D34/100*C34"
$wsAug.Cells.Item(35,1).Value = 30
$wsAug.Cells.Item(35,3).Value = 5
$wsAug.Cells.Item(35,4).Value = 12.6
$wsAug.Cells.Item(35,5).Formula = "=D35/100*C35"

# -- summary block -----------------------------------------------------
$wsAug.Cells.Item(37,2).Value = "sum"
$wsAug.Cells.Item(37,3).Formula = "=SUM(C6:C35)"
$wsAug.Cells.Item(37,5).Formula = "=SUM(E6:E35)"

$wsAug.Cells.Item(39,2).Value = "average kWh/100km"
$wsAug.Cells.Item(39,4).Formula = "=E37*100/C37"

$wsAug.Cells.Item(41,2).Value = "estimated range (with 28kWh)"
$wsAug.Cells.Item(41,5).Formula = "=28/D39*100"
$wsAug.Cells.Item(41,9).Value = "estimated range (with 23kWh)"
$wsAug.Cells.Item(41,12).Formula = "=23/D39*100"

$wsAug.Cells.Item(42,2).Value = "observed GOM value at 100% SOC"
$wsAug.Cells.Item(42,5).Value = 233

$wsAug.Cells.Item(44,2).Value = "Conclusion: good match"

# -- view state: August is the active/visible tab, June no longer is --
$wsAug.Range("E44").Select()
$excel.ActiveWindow.ScrollRow = 25

$wsJune.Range("D45").Select()
$excel.ActiveWindow.ScrollRow = 23

$wsAug.Activate()
$wsAug.Range("E44").Select()
